# "Age of Egypt" edit:
#   1. Remove the "Meta description: ..." paragraph that follows the title.
#   2. Insert a new bold paragraph "Play Age of Egypt Free Playtech Slot
#      Online" right before the closing "Prompt: ..." paragraph.
#   3. Replace the "Prompt: Create a colorful..." text with the meta
#      description text (keeping that paragraph's italic formatting).

$d = $word.ActiveDocument

# --- Step 1: locate & delete the whole "Meta description" paragraph ---
$metaFindRange = $d.Content
$metaFindRange.Find.Execute(
    "Meta description", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0
) | Out-Null
$metaPara = $metaFindRange.Paragraphs(1)
$metaPara.Range.Delete() | Out-Null

# --- Step 2: insert a new bold paragraph right before the "Prompt: ..." one ---
$promptFindRange = $d.Content
$promptFindRange.Find.Execute(
    "Prompt:", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0
) | Out-Null
$promptPara = $promptFindRange.Paragraphs(1)

$insertRange = $d.Range($promptPara.Range.Start, $promptPara.Range.Start)
$insertRange.InsertParagraphBefore() | Out-Null

# Re-find the "Prompt: ..." paragraph (positions shifted after the insert)
# and grab the now-empty paragraph immediately preceding it - that's the
# placeholder paragraph mark we just created.
$promptFindRange2 = $d.Content
$promptFindRange2.Find.Execute(
    "Prompt:", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0
) | Out-Null
$promptPara2 = $promptFindRange2.Paragraphs(1)
$newEmptyPara = $promptPara2.Previous()

# Replace that placeholder paragraph's whole range (start..end, i.e.
# including its pilcrow) with fresh OOXML so the result matches the
# document's "leading empty run" pattern (<w:r/> followed by the formatted
# run) instead of inheriting formatting from the following paragraph.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Age of Egypt Free Playtech Slot Online</w:t></w:r></w:p>' + `
  '<w:p><w:pPr/></w:p>' + `
  '</w:body></w:document>' + `
  '</pkg:xmlData></pkg:part></pkg:package>'
$target = $d.Range($newEmptyPara.Range.Start, $newEmptyPara.Range.End)
$target.InsertXML($xml) | Out-Null

# --- Step 3: replace the "Prompt: ..." text with the new description text,
#     keeping the paragraph's existing (italic) run formatting ---
$d.Content.Find.Execute(
    "Prompt: Create a colorful cartoon-style feature image for the online slot " + [char]34 + "Age of Egypt" + [char]34 + ", featuring a happy Maya warrior with glasses. The image should have a fun and adventurous feel, with the warrior holding a treasure from the game and standing in front of a famous landmark from ancient Egypt, such as the Great Sphinx or the pyramids. Use bright and bold colors to catch the attention of potential players, and include the game title and the Playtech logo as well.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Explore the theme of ancient Egypt with Age of Egypt by Playtech. Play this online slot game for free and search for hidden treasures with bonus features.",
    2
) | Out-Null
